# Add a new column K "intervention_type" with values for each clinical trial row,
# mirroring the header style used by the existing header row (A1:J1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell
$ws.Range("K1").Value = "intervention_type"

# Copy the header formatting (bold font, borders, centered/top alignment)
# from an existing header cell so the new header matches the others exactly.
$ws.Range("A1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for rows 2-15
$values = @{
    2  = "DEVICE"
    3  = "OTHER"
    4  = "PROCEDURE"
    5  = "PROCEDURE"
    6  = "BEHAVIORAL"
    7  = "BEHAVIORAL"
    8  = "DRUG"
    9  = "DEVICE"
    10 = "OTHER"
    11 = "OTHER"
    12 = "DEVICE"
    13 = "OTHER"
    14 = "BEHAVIORAL"
    15 = "PROCEDURE"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 11).Value = $values[$row]
}

# Row 16 has no intervention_type value, but the column still needs a cell
# present there (empty), consistent with other blank cells in that row.
# Copy the (default/unstyled) format from a plain data cell so the new cell
# is created without introducing a new style.
$ws.Range("A2").Copy()
$ws.Range("K16").PasteSpecial(-4122)  # xlPasteFormats
